$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header labels in row 1 (columns A-G); H1/I1/J1 stay the same.
$ws.Range("A1").Value = "codigo"
$ws.Range("B1").Value = "descripcion"
$ws.Range("C1").Value = "codigo_barras"
$ws.Range("D1").Value = "precio_menudeo"
$ws.Range("E1").Value = "precio_mayoreo"
$ws.Range("F1").Value = "precio_caja"
$ws.Range("G1").Value = "precio_socio"

# Column G (rows 2-15) lost its (no-op) applied-number-format style.
$ws.Range("G2:G15").ClearFormats()

# Move the active selection from C6 to F6.
$ws.Range("F6").Select() | Out-Null
